$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 31 new rows before the old row 278 (which will become row 309)
$ws.Range("A270:A300").EntireRow.Insert()

# Block 1: Rambam 10 95% 30deg (rows 270-289)
$ws.Range("A270").Value = 44171
$ws.Range("B270").Value = 'Rambam 10 95% 30deg'
$ws.Range("C270").Value = 301095
$ws.Range("D270").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample1 30deg\Capture 1'
$ws.Range("A271").Value = 44171
$ws.Range("B271").Value = 'Rambam 10 95% 30deg'
$ws.Range("C271").Value = 301095
$ws.Range("D271").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample1 30deg\Capture 2'
$ws.Range("A272").Value = 44171
$ws.Range("B272").Value = 'Rambam 10 95% 30deg'
$ws.Range("C272").Value = 301095
$ws.Range("D272").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample1 30deg\Capture 3'
$ws.Range("A273").Value = 44171
$ws.Range("B273").Value = 'Rambam 10 95% 30deg'
$ws.Range("C273").Value = 301095
$ws.Range("D273").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample1 30deg\Capture 4'
$ws.Range("A274").Value = 44171
$ws.Range("B274").Value = 'Rambam 10 95% 30deg'
$ws.Range("C274").Value = 301095
$ws.Range("D274").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample1 30deg\Capture 5'
$ws.Range("A275").Value = 44171
$ws.Range("B275").Value = 'Rambam 10 95% 30deg'
$ws.Range("C275").Value = 301095
$ws.Range("D275").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample1 30deg\Capture 6'
$ws.Range("A276").Value = 44171
$ws.Range("B276").Value = 'Rambam 10 95% 30deg'
$ws.Range("C276").Value = 301095
$ws.Range("D276").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 1'
$ws.Range("A277").Value = 44171
$ws.Range("B277").Value = 'Rambam 10 95% 30deg'
$ws.Range("C277").Value = 301095
$ws.Range("D277").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 2'
$ws.Range("A278").Value = 44171
$ws.Range("B278").Value = 'Rambam 10 95% 30deg'
$ws.Range("C278").Value = 301095
$ws.Range("D278").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 3'
$ws.Range("A279").Value = 44171
$ws.Range("B279").Value = 'Rambam 10 95% 30deg'
$ws.Range("C279").Value = 301095
$ws.Range("D279").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 4'
$ws.Range("A280").Value = 44171
$ws.Range("B280").Value = 'Rambam 10 95% 30deg'
$ws.Range("C280").Value = 301095
$ws.Range("D280").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 6'
$ws.Range("A281").Value = 44171
$ws.Range("B281").Value = 'Rambam 10 95% 30deg'
$ws.Range("C281").Value = 301095
$ws.Range("D281").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 7'
$ws.Range("A282").Value = 44171
$ws.Range("B282").Value = 'Rambam 10 95% 30deg'
$ws.Range("C282").Value = 301095
$ws.Range("D282").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix1 rambam10 95% sample2 30deg\Capture 8'
$ws.Range("A283").Value = 44171
$ws.Range("B283").Value = 'Rambam 10 95% 30deg'
$ws.Range("C283").Value = 301095
$ws.Range("D283").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 1'
$ws.Range("A284").Value = 44171
$ws.Range("B284").Value = 'Rambam 10 95% 30deg'
$ws.Range("C284").Value = 301095
$ws.Range("D284").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 3'
$ws.Range("A285").Value = 44171
$ws.Range("B285").Value = 'Rambam 10 95% 30deg'
$ws.Range("C285").Value = 301095
$ws.Range("D285").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 5'
$ws.Range("A286").Value = 44171
$ws.Range("B286").Value = 'Rambam 10 95% 30deg'
$ws.Range("C286").Value = 301095
$ws.Range("D286").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 6'
$ws.Range("A287").Value = 44171
$ws.Range("B287").Value = 'Rambam 10 95% 30deg'
$ws.Range("C287").Value = 301095
$ws.Range("D287").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 7'
$ws.Range("A288").Value = 44171
$ws.Range("B288").Value = 'Rambam 10 95% 30deg'
$ws.Range("C288").Value = 301095
$ws.Range("D288").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 8'
$ws.Range("A289").Value = 44171
$ws.Range("B289").Value = 'Rambam 10 95% 30deg'
$ws.Range("C289").Value = 301095
$ws.Range("D289").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample2 30deg\Capture 9'

# Block 2: Rambam 10 95% 25deg (rows 291-295)
$ws.Range("A291").Value = 44171
$ws.Range("B291").Value = 'Rambam 10 95% 25deg'
$ws.Range("C291").Value = 251095
$ws.Range("D291").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample1 25deg\Capture 1'
$ws.Range("A292").Value = 44171
$ws.Range("B292").Value = 'Rambam 10 95% 25deg'
$ws.Range("C292").Value = 251095
$ws.Range("D292").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample1 25deg\Capture 2'
$ws.Range("A293").Value = 44171
$ws.Range("B293").Value = 'Rambam 10 95% 25deg'
$ws.Range("C293").Value = 251095
$ws.Range("D293").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample1 25deg\Capture 3'
$ws.Range("A294").Value = 44171
$ws.Range("B294").Value = 'Rambam 10 95% 25deg'
$ws.Range("C294").Value = 251095
$ws.Range("D294").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample1 25deg\Capture 4'
$ws.Range("A295").Value = 44171
$ws.Range("B295").Value = 'Rambam 10 95% 25deg'
$ws.Range("C295").Value = 251095
$ws.Range("D295").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_06 temperature\mix2 rambam10 95% sample1 25deg\Capture 5'

# Block 3: Rambam 10 95% 30deg 300nM calyculin and 13mM Mg (rows 297-307)
$ws.Range("A297").Value = 44172
$ws.Range("B297").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C297").Value = 30300131095
$ws.Range("D297").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample1 30deg\63x\Capture 2'
$ws.Range("A298").Value = 44172
$ws.Range("B298").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C298").Value = 30300131095
$ws.Range("D298").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample1 30deg\63x\Capture 3'
$ws.Range("A299").Value = 44172
$ws.Range("B299").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C299").Value = 30300131095
$ws.Range("D299").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample1 30deg\63x\Capture 4'
$ws.Range("A300").Value = 44172
$ws.Range("B300").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C300").Value = 30300131095
$ws.Range("D300").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample1 30deg\63x\Capture 5'
$ws.Range("A301").Value = 44172
$ws.Range("B301").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C301").Value = 30300131095
$ws.Range("D301").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample1 30deg\63x\Capture 6'
$ws.Range("A302").Value = 44172
$ws.Range("B302").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C302").Value = 30300131095
$ws.Range("D302").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample1 30deg\63x\Capture 7'
$ws.Range("A303").Value = 44172
$ws.Range("B303").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C303").Value = 30300131095
$ws.Range("D303").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample2 30deg\Capture 1'
$ws.Range("A304").Value = 44172
$ws.Range("B304").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C304").Value = 30300131095
$ws.Range("D304").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample2 30deg\Capture 2'
$ws.Range("A305").Value = 44172
$ws.Range("B305").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C305").Value = 30300131095
$ws.Range("D305").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample2 30deg\Capture 3'
$ws.Range("A306").Value = 44172
$ws.Range("B306").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C306").Value = 30300131095
$ws.Range("D306").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample2 30deg\Capture 4'
$ws.Range("A307").Value = 44172
$ws.Range("B307").Value = 'Rambam 10 95% 30deg 300nM calyculin and 13mM Mg '
$ws.Range("C307").Value = 30300131095
$ws.Range("D307").Value = 'W:\phkinnerets\storage\analysis\Niv\rambam10\95%\2020_12_07 temperature, Mg, and Calyculin\mix1 95% rambam10 with 300nm calyculin A and 13mM Mg sample2 30deg\Capture 5'
